$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L data mirroring column K for rows 3 and 4 (values + formats)
$ws.Range("K3:K4").Copy()
$ws.Range("L3:L4").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("K3:K4").Copy()
$ws.Range("L3:L4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to M12
$ws.Range("M12").Select()
